$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = "Mon Feb 24 23:04:50 EST 2025"
$ws.Range("B3").Value = "Mon Feb 24 23:05:03 EST 2025"
$ws.Range("B4").Value = "Mon Feb 24 23:05:17 EST 2025"
$ws.Range("B5").Value = "Mon Feb 24 23:05:31 EST 2025"
$ws.Range("B6").Value = "Mon Feb 24 23:05:44 EST 2025"
$ws.Range("B7").Value = "Mon Feb 24 23:05:57 EST 2025"
